$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")
$ws.Range("D4").Value = "Incorrect login or password"
